$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the prior-list scenario text in column P (rows 2-33) to the new
#    two-component list used for the simulation study.
$ws.Range("P2:P33").Value = "list(list(shape=1,scale=1),list(shape=1,scale=0.1))"

# 2. Bump sigma_true / mu_sd_cf_prior scenario values for the "large" sigma
#    rows (10-17 and 26-33) from 0.7 to 2.
$ws.Range("G10:H17").Value = 2
$ws.Range("G26:H33").Value = 2

# 3. Widen columns H and P so the longer text fits (matches Excel's bestFit).
$ws.Columns.Item(8).ColumnWidth = 13.87
$ws.Columns.Item(16).ColumnWidth = 45.45

# 4. Restore the selection to the range that was just edited.
$ws.Range("P2:P33").Select() | Out-Null
